$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13: aa / Agree/Accept -> ba / Appreciation
$ws.Range("I13").Value = "ba"
$ws.Range("J13").Value = "Appreciation"

# Row 31: sv / Statement-opinion -> sd / Statement-non-opinion
$ws.Range("I31").Value = "sd"
$ws.Range("J31").Value = "Statement-non-opinion"

# Row 45: sd / Statement-non-opinion -> ba / Appreciation
$ws.Range("I45").Value = "ba"
$ws.Range("J45").Value = "Appreciation"

# Row 46: qy / Yes-No-Question -> ba / Appreciation
$ws.Range("I46").Value = "ba"
$ws.Range("J46").Value = "Appreciation"
